# Populate new flight-arrival rows 76-92 (Tuesday, Jan 10) on the "Main Data" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(76, 1).Value = 75.0
$ws.Range('B76').Value = 'Tuesday, Jan 10'
$ws.Range('C76').Value = '12:05 AM'
$ws.Range('D76').Value = 'FR1969'
$ws.Range('E76').Value = 'Madrid'
$ws.Range('F76').Value = '(MAD)'
$ws.Range('G76').Value = 'Ryanair '
$ws.Range('H76').Value = 'B38M'
$ws.Range('I76').Value = '(SP-RZO)'
$ws.Range('J76').Value = '12:07 AM'
$ws.Range('K76').Borders.LineStyle = 0
$ws.Range('L76').Value = '0 hours, 2 minutes'
$ws.Range('M76').Borders.LineStyle = 0

$ws.Cells.Item(77, 1).Value = 76.0
$ws.Range('B77').Value = 'Tuesday, Jan 10'
$ws.Range('C77').Value = '9:05 AM'
$ws.Range('D77').Value = 'FR1923'
$ws.Range('E77').Value = 'Budapest'
$ws.Range('F77').Value = '(BUD)'
$ws.Range('G77').Value = 'Ryanair '
$ws.Range('H77').Value = 'B738'
$ws.Range('I77').Value = '(SP-RSC)'
$ws.Range('J77').Value = '8:21 AM'
$ws.Range('K77').Borders.LineStyle = 0
$ws.Range('L77').Value = '0 hours, -44 minutes'
$ws.Range('M77').Borders.LineStyle = 0

$ws.Cells.Item(78, 1).Value = 77.0
$ws.Range('B78').Value = 'Tuesday, Jan 10'
$ws.Range('C78').Value = '10:05 AM'
$ws.Range('D78').Value = 'FR8011'
$ws.Range('E78').Value = 'Dublin'
$ws.Range('F78').Value = '(DUB)'
$ws.Range('G78').Value = 'Ryanair '
$ws.Range('H78').Value = 'B38M'
$ws.Range('I78').Value = '(EI-HGX)'
$ws.Range('J78').Value = '9:50 AM'
$ws.Range('K78').Borders.LineStyle = 0
$ws.Range('L78').Value = '0 hours, -15 minutes'
$ws.Range('M78').Borders.LineStyle = 0

$ws.Cells.Item(79, 1).Value = 78.0
$ws.Range('B79').Value = 'Tuesday, Jan 10'
$ws.Range('C79').Value = '10:50 AM'
$ws.Range('D79').Value = 'FR1967'
$ws.Range('E79').Value = 'Milan'
$ws.Range('F79').Value = '(BGY)'
$ws.Range('G79').Value = 'Buzz '
$ws.Range('H79').Value = 'B38M'
$ws.Range('I79').Value = '(SP-RZG)'
$ws.Range('J79').Value = '10:52 AM'
$ws.Range('K79').Borders.LineStyle = 0
$ws.Range('L79').Value = '0 hours, 2 minutes'
$ws.Range('M79').Borders.LineStyle = 0

$ws.Cells.Item(80, 1).Value = 79.0
$ws.Range('B80').Value = 'Tuesday, Jan 10'
$ws.Range('C80').Value = '11:10 AM'
$ws.Range('D80').Value = 'FR2019'
$ws.Range('E80').Value = 'London'
$ws.Range('F80').Value = '(STN)'
$ws.Range('G80').Value = 'Ryanair '
$ws.Range('H80').Value = 'B738'
$ws.Range('I80').Value = '(SP-RKL)'
$ws.Range('J80').Value = '11:15 AM'
$ws.Range('K80').Borders.LineStyle = 0
$ws.Range('L80').Value = '0 hours, 5 minutes'
$ws.Range('M80').Borders.LineStyle = 0

$ws.Cells.Item(81, 1).Value = 80.0
$ws.Range('B81').Value = 'Tuesday, Jan 10'
$ws.Range('C81').Value = '11:45 AM'
$ws.Range('D81').Value = 'FR1941'
$ws.Range('E81').Value = 'Bristol'
$ws.Range('F81').Value = '(BRS)'
$ws.Range('G81').Value = 'Ryanair '
$ws.Range('H81').Value = 'B738'
$ws.Range('I81').Value = '(SP-RKT)'
$ws.Range('J81').Value = '11:43 AM'
$ws.Range('K81').Borders.LineStyle = 0
$ws.Range('L81').Value = '0 hours, -2 minutes'
$ws.Range('M81').Borders.LineStyle = 0

$ws.Cells.Item(82, 1).Value = 81.0
$ws.Range('B82').Value = 'Tuesday, Jan 10'
$ws.Range('C82').Value = '11:55 AM'
$ws.Range('D82').Value = 'FR1889'
$ws.Range('E82').Value = 'Paris'
$ws.Range('F82').Value = '(BVA)'
$ws.Range('G82').Value = 'Ryanair '
$ws.Range('H82').Value = 'B38M'
$ws.Range('I82').Value = '(SP-RZO)'
$ws.Range('J82').Value = '11:54 AM'
$ws.Range('K82').Borders.LineStyle = 0
$ws.Range('L82').Value = '0 hours, -1 minutes'
$ws.Range('M82').Borders.LineStyle = 0

$ws.Cells.Item(83, 1).Value = 82.0
$ws.Range('B83').Value = 'Tuesday, Jan 10'
$ws.Range('C83').Value = '12:05 PM'
$ws.Range('D83').Value = 'FR1112'
$ws.Range('E83').Value = 'Rome'
$ws.Range('F83').Value = '(CIA)'
$ws.Range('G83').Value = 'Ryanair '
$ws.Range('H83').Value = 'B738'
$ws.Range('I83').Value = '(SP-RKP)'
$ws.Range('J83').Value = '11:48 AM'
$ws.Range('K83').Borders.LineStyle = 0
$ws.Range('L83').Value = '0 hours, -17 minutes'
$ws.Range('M83').Borders.LineStyle = 0

$ws.Cells.Item(84, 1).Value = 83.0
$ws.Range('B84').Value = 'Tuesday, Jan 10'
$ws.Range('C84').Value = '2:10 PM'
$ws.Range('D84').Value = 'FR4178'
$ws.Range('E84').Value = 'Pafos'
$ws.Range('F84').Value = '(PFO)'
$ws.Range('G84').Value = 'Ryanair '
$ws.Range('H84').Value = 'B38M'
$ws.Range('I84').Value = '(SP-RZI)'
$ws.Range('J84').Value = '2:02 PM'
$ws.Range('K84').Borders.LineStyle = 0
$ws.Range('L84').Value = '0 hours, -8 minutes'
$ws.Range('M84').Borders.LineStyle = 0

$ws.Cells.Item(85, 1).Value = 84.0
$ws.Range('B85').Value = 'Tuesday, Jan 10'
$ws.Range('C85').Value = '5:35 PM'
$ws.Range('D85').Value = 'FR2669'
$ws.Range('E85').Value = 'London'
$ws.Range('F85').Value = '(STN)'
$ws.Range('G85').Value = 'Ryanair '
$ws.Range('H85').Value = 'B738'
$ws.Range('I85').Value = '(SP-RKT)'
$ws.Range('J85').Value = '6:26 PM'
$ws.Range('K85').Borders.LineStyle = 0
$ws.Range('L85').Value = '0 hours, 51 minutes'
$ws.Range('M85').Borders.LineStyle = 0

$ws.Cells.Item(86, 1).Value = 85.0
$ws.Range('B86').Value = 'Tuesday, Jan 10'
$ws.Range('C86').Value = '6:00 PM'
$ws.Range('D86').Value = 'FR1969'
$ws.Range('E86').Value = 'Madrid'
$ws.Range('F86').Value = '(MAD)'
$ws.Range('G86').Value = 'Ryanair '
$ws.Range('H86').Value = 'B738'
$ws.Range('I86').Value = '(EI-ENL)'
$ws.Range('J86').Value = '6:42 PM'
$ws.Range('K86').Borders.LineStyle = 0
$ws.Range('L86').Value = '0 hours, 42 minutes'
$ws.Range('M86').Borders.LineStyle = 0

$ws.Cells.Item(87, 1).Value = 86.0
$ws.Range('B87').Value = 'Tuesday, Jan 10'
$ws.Range('C87').Value = '6:20 PM'
$ws.Range('D87').Value = 'FR1939'
$ws.Range('E87').Value = 'Gothenburg'
$ws.Range('F87').Value = '(GOT)'
$ws.Range('G87').Value = 'Ryanair '
$ws.Range('H87').Value = 'B38M'
$ws.Range('I87').Value = '(SP-RZI)'
$ws.Range('J87').Value = '6:11 PM'
$ws.Range('K87').Borders.LineStyle = 0
$ws.Range('L87').Value = '0 hours, -9 minutes'
$ws.Range('M87').Borders.LineStyle = 0

$ws.Cells.Item(88, 1).Value = 87.0
$ws.Range('B88').Value = 'Tuesday, Jan 10'
$ws.Range('C88').Value = '8:10 PM'
$ws.Range('D88').Value = 'FR2264'
$ws.Range('E88').Value = 'Lisbon'
$ws.Range('F88').Value = '(LIS)'
$ws.Range('G88').Value = 'Ryanair '
$ws.Range('H88').Value = 'B738'
$ws.Range('I88').Value = '(SP-RSP)'
$ws.Range('J88').Value = '8:26 PM'
$ws.Range('K88').Borders.LineStyle = 0
$ws.Range('L88').Value = '0 hours, 16 minutes'
$ws.Range('M88').Borders.LineStyle = 0

$ws.Cells.Item(89, 1).Value = 88.0
$ws.Range('B89').Value = 'Tuesday, Jan 10'
$ws.Range('C89').Value = '9:00 PM'
$ws.Range('D89').Value = 'FR4059'
$ws.Range('E89').Value = 'Malaga'
$ws.Range('F89').Value = '(AGP)'
$ws.Range('G89').Value = 'Buzz '
$ws.Range('H89').Value = 'B38M'
$ws.Range('I89').Value = '(SP-RZG)'
$ws.Range('J89').Value = '9:02 PM'
$ws.Range('K89').Borders.LineStyle = 0
$ws.Range('L89').Value = '0 hours, 2 minutes'
$ws.Range('M89').Borders.LineStyle = 0

$ws.Cells.Item(90, 1).Value = 89.0
$ws.Range('B90').Value = 'Tuesday, Jan 10'
$ws.Range('C90').Value = '10:10 PM'
$ws.Range('D90').Value = 'FR1573'
$ws.Range('E90').Value = 'Vienna'
$ws.Range('F90').Value = '(VIE)'
$ws.Range('G90').Value = 'Ryanair '
$ws.Range('H90').Value = 'B38M'
$ws.Range('I90').Value = '(SP-RZI)'
$ws.Range('J90').Value = '10:13 PM'
$ws.Range('K90').Borders.LineStyle = 0
$ws.Range('L90').Value = '0 hours, 3 minutes'
$ws.Range('M90').Borders.LineStyle = 0

$ws.Cells.Item(91, 1).Value = 90.0
$ws.Range('B91').Value = 'Tuesday, Jan 10'
$ws.Range('C91').Value = '11:25 PM'
$ws.Range('D91').Value = 'FR1907'
$ws.Range('E91').Value = 'Milan'
$ws.Range('F91').Value = '(BGY)'
$ws.Range('G91').Value = 'Ryanair '
$ws.Range('H91').Value = 'B738'
$ws.Range('I91').Value = '(SP-RKP)'
$ws.Range('J91').Value = 'Diverted to WAW'
$ws.Range('K91').Borders.LineStyle = 0
$ws.Range('L91').Borders.LineStyle = 0
$ws.Range('M91').Borders.LineStyle = 0

$ws.Cells.Item(92, 1).Value = 91.0
$ws.Range('B92').Value = 'Tuesday, Jan 10'
$ws.Range('C92').Value = '11:45 PM'
$ws.Range('D92').Value = 'FR1021'
$ws.Range('E92').Value = 'London'
$ws.Range('F92').Value = '(STN)'
$ws.Range('G92').Value = 'Ryanair '
$ws.Range('H92').Value = 'B38M'
$ws.Range('I92').Value = '(SP-RZO)'
$ws.Range('J92').Value = '12:11 AM'
$ws.Range('K92').Borders.LineStyle = 0
$ws.Range('L92').Value = '0 hours, 26 minutes'
$ws.Range('M92').Borders.LineStyle = 0
